# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.353.74"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "3.467.79"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.76"
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.85"
$ws.Range("E6").Value = "  +2.98%  "
$ws.Range("E7").Value = "  +3.79%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "3.465.42"
$ws.Range("E9").Value = "  -1.06%  "
$ws.Range("E10").Value = "  +6.79%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.98"
$ws.Range("E11").Value = "  -2.24%  "
$ws.Range("E12").Value = "  -0.40%  "
$ws.Range("D13").Value = "4.065.14"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.06"
$ws.Range("E14").Value = "  +2.90%  "
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "67.347.60"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "3.469.84"
$ws.Range("E18").Value = "  -1.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  -1.63%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.12"
$ws.Range("E20").Value = "  -3.39%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "394.29"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.93"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  -0.02%  "
$ws.Range("E25").Value = "  -0.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.73"
$ws.Range("E26").Value = "  -2.44%  "
$ws.Range("E27").Value = "  -0.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.36"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E29").Value = "  -2.68%  "
$ws.Range("E30").Value = "  +0.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.12"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.40"
$ws.Range("E32").Value = "  -1.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.06"
$ws.Range("E33").Value = "  -0.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "23.51"
$ws.Range("E34").Value = "  -0.72%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "160.79"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("E40").Value = "  +11.94%  "
$ws.Range("E41").Value = "  -3.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.66"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.72"
$ws.Range("E43").Value = "  -4.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "26.14"
$ws.Range("E44").Value = "  -1.35%  "
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("D46").Value = "2.759.84"
$ws.Range("E46").Value = "  -1.95%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "26.29"
$ws.Range("E47").Value = "  -3.18%  "
$ws.Range("E48").Value = "  -2.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0298"
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "326.24"
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.04"
$ws.Range("E51").Value = "  -3.17%  "

Write-Host "Updated cryptos list"
